# --- Excel COM-interop script implementing the commit's changes ---
# 1. Rename header cells on the two existing sheets.
# 2. Add a new "PO Forecast" worksheet after "Monthly Trend".
# 3. Populate the new sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper),
#    reusing the same header/date formatting already used on the other sheets.

$wb = $excel.ActiveWorkbook

# --- 1. Update existing headers -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" worksheet at the end of the workbook ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used by the other sheets (0.75in/1in/0.5in => 54pt/72pt/36pt)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- 3. Populate headers (copy formatting from an existing header cell) --------
$wsWeekly.Range("A1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Apply the date formatting used by the other sheets' date column to column A ---
$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 5. Populate data rows ---------------------------------------------------------
$data = @(
    @(45109.99999999999, 58, -7.471373698781782, 120.0010360288059),
    @(45137.99999999999, 76, 18.27239335141422, 139.5356423222848),
    @(45144.99999999999, 80, 19.83210755150064, 143.4819621566729),
    @(45151.99999999999, 85, 19.21260856527844, 145.7042797656568),
    @(45165.99999999999, 94, 29.22218796395838, 157.0494951541563),
    @(45172.99999999999, 99, 33.24592925773462, 166.6302009270561),
    @(45179.99999999999, 103, 39.72634262622087, 167.3633371461696),
    @(45186.99999999999, 108, 44.63930551408504, 164.9061180725879),
    @(45200.99999999999, 117, 53.81620390735185, 178.7313055382743),
    @(45207.99999999999, 122, 60.00863366266202, 183.071291685877),
    @(45214.99999999999, 126, 67.74555383378079, 188.0541672850416),
    @(45221.99999999999, 131, 64.58204474776007, 196.7897249429046),
    @(45228.99999999999, 135, 77.07038953291773, 195.9389074627175),
    @(45235.99999999999, 140, 81.00464924137903, 201.5923602813319),
    @(45242.99999999999, 144, 80.09469243900259, 205.4136845945299),
    @(45249.99999999999, 149, 85.99199941362258, 212.4460002981291),
    @(45256.99999999999, 154, 92.93441747617238, 214.9167885554417)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
